$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for the crypto rows that changed.
# Values are prefixed with a literal apostrophe so Excel stores them as
# text (matching the original inline-string cells) instead of coercing
# numeric-looking strings (e.g. '96.69') into numbers.

$ws.Range('D2').Value = "'42.958.46"
$ws.Range('E2').Value = "'  -1.36%  "
$ws.Range('D3').Value = "'2.578.34"
$ws.Range('E3').Value = "'  -0.47%  "
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = "'  -0.08%  "
$ws.Range('D5').Value = "'302.94"
$ws.Range('E5').Value = "'  +0.65%  "
$ws.Range('D6').Value = "'96.69"
$ws.Range('E6').Value = "'  +0.34%  "
$ws.Range('D7').Value = "'0.576"
$ws.Range('E7').Value = "'  -0.65%  "
$ws.Range('E8').Value = "'  -0.09%  "
$ws.Range('D9').Value = "'0.549"
$ws.Range('E9').Value = "'  -1.60%  "
$ws.Range('D10').Value = "'36.72"
$ws.Range('E10').Value = "'  +0.03%  "
$ws.Range('D11').Value = "'0.0814"
$ws.Range('E11').Value = "'  -0.33%  "
$ws.Range('D12').Value = "'7.65"
$ws.Range('E12').Value = "'  -1.94%  "
$ws.Range('E13').Value = "'  +6.42%  "
$ws.Range('D14').Value = "'2.601.83"
$ws.Range('E14').Value = "'  +0.45%  "
$ws.Range('D15').Value = "'0.886"
$ws.Range('E15').Value = "'  -0.47%  "
$ws.Range('D16').Value = "'14.29"
$ws.Range('E16').Value = "'  -0.23%  "
$ws.Range('D17').Value = "'43.012.07"
$ws.Range('E17').Value = "'  -1.18%  "
$ws.Range('E18').Value = "'  +5.42%  "
$ws.Range('D19').Value = "'0.0₃0999"
$ws.Range('E19').Value = "'  +1.95%  "
$ws.Range('D20').Value = "'6.67"
$ws.Range('E20').Value = "'  -0.05%  "
$ws.Range('D21').Value = "'71.95"
$ws.Range('E21').Value = "'  -1.13%  "
$ws.Range('E22').Value = "'  -4.06%  "
$ws.Range('E23').Value = "'  +0.89%  "
$ws.Range('E24').Value = "'  -3.87%  "
$ws.Range('D25').Value = "'28.97"
$ws.Range('E25').Value = "'  -0.65%  "
$ws.Range('E26').Value = "'  -0.14%  "
$ws.Range('E27').Value = "'  +0.53%  "
$ws.Range('D28').Value = "'37.65"
$ws.Range('E28').Value = "'  -0.57%  "
$ws.Range('E29').Value = "'  -2.57%  "
$ws.Range('E30').Value = "'  -0.77%  "
$ws.Range('D31').Value = "'155.00"
$ws.Range('E31').Value = "'  +1.90%  "
$ws.Range('E32').Value = "'  -1.06%  "
$ws.Range('D33').Value = "'3.43"
$ws.Range('E33').Value = "'  -5.13%  "
$ws.Range('E34').Value = "'  -1.33%  "
$ws.Range('E35').Value = "'  -0.83%  "
$ws.Range('D36').Value = "'18.15"
$ws.Range('E36').Value = "'  +9.05%  "
$ws.Range('E37').Value = "'  -2.98%  "
$ws.Range('E38').Value = "'  -0.39%  "
$ws.Range('D39').Value = "'23.04"
$ws.Range('E39').Value = "'  -5.34%  "
$ws.Range('D40').Value = "'2.25"
$ws.Range('E40').Value = "'  +40.86%  "
$ws.Range('D41').Value = "'3.45"
$ws.Range('E41').Value = "'  -4.91%  "
$ws.Range('E43').Value = "'  +0.94%  "
$ws.Range('D44').Value = "'2.086.21"
$ws.Range('E44').Value = "'  +2.04%  "
$ws.Range('D45').Value = "'0.998"
$ws.Range('E45').Value = "'  -0.03%  "
$ws.Range('D46').Value = "'9.25"
$ws.Range('E46').Value = "'  +1.96%  "
$ws.Range('D47').Value = "'85.38"
$ws.Range('E47').Value = "'  -2.94%  "
$ws.Range('D48').Value = "'76.83"
$ws.Range('E48').Value = "'  +11.18%  "
$ws.Range('D49').Value = "'106.71"
$ws.Range('E49').Value = "'  +0.88%  "
$ws.Range('D50').Value = "'2.828.03"
$ws.Range('E50').Value = "'  -0.75%  "
$ws.Range('D51').Value = "'1.69"
$ws.Range('E51').Value = "'  +2.74%  "

# Reset style on the touched range so the quote-prefix flag Excel adds
# for text-forced numeric values does not linger on the cells (the
# source cells carry no explicit style).
$ws.Range("D2:E51").Style = "Normal"

